$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 16-27: columns C (Doc #), D (Name), E (Period)
# Grouped by worker (JHONNY then MAURICIO), periods sorted descending (1709 -> 1704)

$docJhonny = "1050953042"
$nameJhonny = "JHONNY JAVIER CARRILLO GARCIA"

$docMauricio = "1047378133"
$nameMauricio = "MAURICIO JAVIER TORRES ANGULO"

$periods = @("1709", "1708", "1707", "1706", "1705", "1704")

# Rows 16-21 : JHONNY, periods 1709..1704
$row = 16
foreach ($p in $periods) {
    $ws.Cells.Item($row, 3).Value = $docJhonny
    $ws.Cells.Item($row, 4).Value = $nameJhonny
    $ws.Cells.Item($row, 5).Value = $p
    $row = $row + 1
}

# Rows 22-27 : MAURICIO, periods 1709..1704
foreach ($p in $periods) {
    $ws.Cells.Item($row, 3).Value = $docMauricio
    $ws.Cells.Item($row, 4).Value = $nameMauricio
    $ws.Cells.Item($row, 5).Value = $p
    $row = $row + 1
}
